# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: "Created functions to get season record" — the
# team's W/L/T totals are appended as three new columns (AD, AE, AF)
# and repeated for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — new labels, styled like the existing header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold font + border + alignment)
# from an existing header cell onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values, repeated for every player row (2 through 44).
$wins = 71
$losses = 91
$ties = 0

for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
